# Adding CPE 621 to title slide
#
# Slide 1's subtitle placeholder currently reads:
#   "Final Presentation"
#   "Caleb Stewart & John Wilkes"
# A new line, "CPE 621", is appended as its own paragraph. After that
# edit, the "Caleb Stewart & John Wilkes" line ends up split across two
# runs ("Caleb Stewart & John " + "Wilkes"), which is what PowerPoint
# leaves behind once that text box has been edited/retyped in the UI.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)          # "Subtitle 2" placeholder
$tr  = $shp.TextFrame.TextRange

# --- Append the new "CPE 621" paragraph --------------------------------
$originalText = $tr.Text
$tr.Text = $originalText + "`rCPE 621"

# --- Re-split "Caleb Stewart & John Wilkes" into two runs --------------
# Locate the "Caleb Stewart & John Wilkes" line within the text and split
# off the "Caleb Stewart & John " portion into its own run, leaving
# "Wilkes" as a separate trailing run.
$namesLine  = "Caleb Stewart & John Wilkes"
$firstPart  = "Caleb Stewart & John "

$lineOffset = $originalText.IndexOf($namesLine)        # 0-based offset
$lineStart  = $lineOffset + 1                          # 1-based start

$firstRun = $tr.Characters($lineStart, $firstPart.Length)
$firstRun.Text = $firstPart
